# Regenerate save_data: column G ("K" = strikeouts) is recalculated from the
# updated source (K instead of Strike#). Write the new per-game strikeout
# counts into column G (rows 2-19), leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 3
    13 = 2
    14 = 0
    16 = 2
    17 = 2
    18 = 1
    19 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
